$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, copying the formatting from an existing header cell (H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the data values for columns I and J, rows 2-13
$values = @{
    2  = @(8, 9)
    3  = @(8, 8)
    4  = @(9, 9)
    5  = @(4, 4)
    6  = @(5, 6)
    7  = @(9, 9)
    8  = @(9, 9)
    9  = @(9, 9)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(2, 2)
    13 = @(8, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
